$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "have gained traction..." sentence punctuation fix
# ------------------------------------------------------------------
$r = $d.Content
$old1 = "have gained traction in the science community - in order for data to be FAIR - they have to be"
$new1 = "have gained traction in the science community. In order for data to be FAIR, they have to be"
$null = $r.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ------------------------------------------------------------------
# 2) Capture the text of the existing TDWG / "Darwin Core is a standard"
#    footnote (currently the 17th footnote in the document) so that we
#    can reuse its content for the new footnote inserted after
#    "Thousands of Darwin Core Archives".
# ------------------------------------------------------------------
$dwcNoteText = $d.Footnotes.Item(17).Range.Text

# ------------------------------------------------------------------
# 3) Expand "Thousands of Darwin Core Archives are published" into
#    "Thousands of Darwin Core Archives (DwC-A) containing valuable
#    biodiversity data are published" (leaving room for the new
#    footnote reference to be inserted right after "Archives").
# ------------------------------------------------------------------
$r = $d.Content
$old3 = "Thousands of Darwin Core Archives are published"
$new3 = "Thousands of Darwin Core Archives  (DwC-A) containing valuable biodiversity data are published"
$null = $r.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# ------------------------------------------------------------------
# 4) Insert the new footnote reference right after
#    "Thousands of Darwin Core Archives " using the captured text.
#    This becomes the document's 16th footnote, pushing the former
#    16th/17th footnotes (Poelen / TDWG) down to 17th/18th.
# ------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Thousands of Darwin Core Archives ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$null = $d.Footnotes.Add($r, "", $dwcNoteText)

# ------------------------------------------------------------------
# 5) Remove the now-redundant old TDWG footnote (it has been pushed to
#    position 18 by the insertion above) together with its reference
#    mark in the body text.
# ------------------------------------------------------------------
$d.Footnotes.Item(18).Delete()

# ------------------------------------------------------------------
# 6) Simplify the "Interoperability" sentence, dropping the inline
#    "(TDWG's Darwin Core Archive <footnote>)" aside in favor of a
#    plain mention of "DwC-A" (the footnote reference that used to sit
#    here was already removed in step 5, leaving a double space behind
#    which this replacement also cleans up).
# ------------------------------------------------------------------
$r = $d.Content
$old6 = "nterability through their adoption on a recognized standards for biodiversity data (TDWG" + [char]0x2019 + "s Darwin Core Archive  ), and was able to"
$new6 = "nterability through their adoption on a recognized standard, DwC-A, and was able to"
$null = $r.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)

# ------------------------------------------------------------------
# 7) Typo fix: "keepinging" -> "keeping"
# ------------------------------------------------------------------
$r = $d.Content
$old7 = "euse the archive by keepinging versioned copies as proof of registration."
$new7 = "euse the archive by keeping versioned copies as proof of registration."
$null = $r.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)

# ------------------------------------------------------------------
# 8) "easy" -> "easier"
# ------------------------------------------------------------------
$r = $d.Content
$old8 = "To make it easy to see whether an archive is FAIR according to the methods describe above, you can get your FAIR assessment badge using:"
$new8 = "To make it easier to see whether an archive is FAIR according to the methods describe above, you can get your FAIR assessment badge using:"
$null = $r.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)

# ------------------------------------------------------------------
# 9) "tracked darwin core archive" -> "tracked DwC-A"
# ------------------------------------------------------------------
$r = $d.Content
$old9 = "If an archive reference (by location, uuid, doi) is associated with a tracked darwin core archive, a download badge is generated"
$new9 = "If an archive reference (by location, uuid, doi) is associated with a tracked DwC-A, a download badge is generated"
$null = $r.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
